$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Add new row 6 data: id=5, itemType="Item_Key_A", description="扉のカギ"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Item_Key_A"
$ws.Range("C6").Value = "扉のカギ"

# Update selection to match target state
$ws.Range("K14").Select()
